$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 331, shifting existing rows 331-360 down to 332-361
$ws.Rows("331:331").Insert()

# Populate the new row 331 with its data
$ws.Range("A331").Value = 5
$ws.Range("B331").Value = "Macroferia Regional de Talca"
$ws.Range("C331").Value = "Maule"
$ws.Range("D331").Value = 45166
$ws.Range("E331").Value = 7
$ws.Range("F331").Value = 100112017
$ws.Range("G331").Value = "Apio"
$ws.Range("H331").Value = "Americana (o)"
$ws.Range("I331").Value = "Primera"
$ws.Range("J331").Value = 500
$ws.Range("K331").Value = 6000
$ws.Range("L331").Value = 6000
$ws.Range("M331").Value = 6000
$ws.Range("N331").Value = "`$/docena de matas"
$ws.Range("O331").Value = "Provincia de Limarí"
$ws.Range("P331").Value = 1000
$ws.Range("Q331").Value = 6
$ws.Range("R331").Value = "Hortaliza"
